$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 12
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null

# Row 62
$ws.Range("H62").Value = 4774.647
$ws.Range("I62").Value = 2623.9092
$ws.Range("J62").Value = 8717.666999999999
$ws.Range("K62").Value = 2623.9092
$ws.Range("L62").Value = 8717.666999999999
$ws.Range("M62").Value = -1999.9092
$ws.Range("N62").Value = -9965.666999999999

# Row 65
$ws.Range("H65").Value = 4774.647
$ws.Range("I65").Value = 2623.9092
$ws.Range("J65").Value = 8717.666999999999
$ws.Range("K65").Value = 13119.546
$ws.Range("L65").Value = 43588.335
$ws.Range("M65").Value = -9999.546
$ws.Range("N65").Value = -49828.335

# Row 94
$ws.Range("H94").Value = 5163
$ws.Range("I94").Value = 1023.5
$ws.Range("K94").Value = 1023.5
$ws.Range("M94").Value = -572.5

# Row 112
$ws.Range("H112").Value = 4043.1
$ws.Range("J112").Value = 4043.1
$ws.Range("L112").Value = 12129.3
$ws.Range("N112").Value = -14345.3

# Row 132
$ws.Range("H132").Value = 1088.2
$ws.Range("I132").Value = 1088.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3264.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -734.6000000000004

# Row 135
$ws.Range("H135").Value = 2273189
$ws.Range("I135").Value = 2273189
$ws.Range("K135").Value = 20458701
$ws.Range("M135").Value = -20456166

# Row 137
$ws.Range("H137").Value = 7108.077
$ws.Range("I137").Value = 8092.5
$ws.Range("J137").Value = 5533
$ws.Range("K137").Value = 24277.5
$ws.Range("L137").Value = 16599
$ws.Range("M137").Value = -21727.5
$ws.Range("N137").Value = -21699

# Row 138
$ws.Range("H138").Value = 5163.241
$ws.Range("J138").Value = 5412.9326
$ws.Range("L138").Value = 16238.7978
$ws.Range("N138").Value = -26518.7978

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 5
$ws.Range("H5").Value = 130
$ws.Range("I5").Value = 74
$ws.Range("K5").Value = 74
$ws.Range("M5").Value = 38

# Row 32
$ws.Range("H32").Value = 4161.86
$ws.Range("I32").Value = 4161.86
$ws.Range("K32").Value = 4161.86
$ws.Range("M32").Value = -3874.86

# Row 61
$ws.Range("H61").Value = 55566460
$ws.Range("I61").Value = 7045.6665
$ws.Range("J61").Value = 83346170
$ws.Range("K61").Value = 7045.6665
$ws.Range("L61").Value = 83346170
$ws.Range("M61").Value = -6833.6665
$ws.Range("N61").Value = -83346594

# Row 74
$ws.Range("H74").Value = 4947.0557
$ws.Range("J74").Value = 5795
$ws.Range("L74").Value = 5795
$ws.Range("N74").Value = -7543

# Row 77
$ws.Range("H77").Value = 4947.0557
$ws.Range("J77").Value = 5795
$ws.Range("L77").Value = 28975
$ws.Range("N77").Value = -37711

# Row 122
$ws.Range("H122").Value = 4350
$ws.Range("I122").Value = 4619.9443
$ws.Range("K122").Value = 13859.8329
$ws.Range("M122").Value = -11409.8329

# Row 132
$ws.Range("H132").Value = 4422.346
$ws.Range("I132").Value = 1923.6765
$ws.Range("J132").Value = 9142.056
$ws.Range("K132").Value = 5771.029500000001
$ws.Range("L132").Value = 27426.168
$ws.Range("M132").Value = -3241.029500000001
$ws.Range("N132").Value = -32486.168

# Row 136
$ws.Range("H136").Value = 55566460
$ws.Range("I136").Value = 7045.6665
$ws.Range("J136").Value = 83346170
$ws.Range("K136").Value = 21136.9995
$ws.Range("L136").Value = 250038510
$ws.Range("M136").Value = -18586.9995
$ws.Range("N136").Value = -250043610

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 4
$ws.Range("H4").Value = 130
$ws.Range("I4").Value = 74
$ws.Range("K4").Value = 74
$ws.Range("M4").Value = 41

# Row 20
$ws.Range("H20").Value = 8340463.5
$ws.Range("I20").Value = 10421580
$ws.Range("K20").Value = 10421580
$ws.Range("M20").Value = -10421333

# Row 134
$ws.Range("H134").Value = 5957783
$ws.Range("I134").Value = 9261194
$ws.Range("J134").Value = 11642.934
$ws.Range("K134").Value = 27783582
$ws.Range("L134").Value = 34928.802
$ws.Range("M134").Value = -27781047
$ws.Range("N134").Value = -39998.802

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 6312.523
$ws.Range("J31").Value = 11073.4
$ws.Range("L31").Value = 11073.4
$ws.Range("N31").Value = -11663.4

# Row 34
$ws.Range("H34").Value = 6312.523
$ws.Range("J34").Value = 11073.4
$ws.Range("L34").Value = 11073.4
$ws.Range("N34").Value = -11477.4

# Row 62
$ws.Range("H62").Value = 2995.6667
$ws.Range("I62").Value = 2996.5
$ws.Range("K62").Value = 2996.5
$ws.Range("M62").Value = -2372.5

# Row 65
$ws.Range("H65").Value = 2995.6667
$ws.Range("I65").Value = 2996.5
$ws.Range("K65").Value = 14982.5
$ws.Range("M65").Value = -11862.5

# Row 76
$ws.Range("H76").Value = 5449.875
$ws.Range("I76").Value = 5449.875
$ws.Range("K76").Value = 5449.875
$ws.Range("M76").Value = -5134.875

# Row 79
$ws.Range("H79").Value = 5449.875
$ws.Range("I79").Value = 5449.875
$ws.Range("K79").Value = 5449.875
$ws.Range("M79").Value = -4357.875

# Row 132
$ws.Range("H132").Value = 4172.2207
$ws.Range("I132").Value = 2645.9592
$ws.Range("K132").Value = 7937.8776
$ws.Range("M132").Value = -5407.8776

# Row 134
$ws.Range("H134").Value = 6156.5938
$ws.Range("I134").Value = 2318.3333
$ws.Range("K134").Value = 6954.999899999999
$ws.Range("M134").Value = -4419.999899999999

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 7
$ws.Range("H7").Value = 429.33334
$ws.Range("J7").Value = 938
$ws.Range("L7").Value = 2814
$ws.Range("N7").Value = -3038

# Row 34
$ws.Range("H34").Value = 18000.334
$ws.Range("J34").Value = 18000.334
$ws.Range("L34").Value = 54001.00199999999
$ws.Range("N34").Value = -54169.00199999999

# Row 39
$ws.Range("H39").Value = 9721.777
$ws.Range("J39").Value = 16499
$ws.Range("L39").Value = 49497
$ws.Range("N39").Value = -50085

# Row 107
$ws.Range("H107").Value = 5238921
$ws.Range("J107").Value = 5527177.5
$ws.Range("L107").Value = 16581532.5
$ws.Range("N107").Value = -16585372.5

# Row 122
$ws.Range("H122").Value = 1741617.9
$ws.Range("I122").Value = 2572918.5
$ws.Range("J122").Value = 910317.25
$ws.Range("K122").Value = 23156266.5
$ws.Range("L122").Value = 8192855.25
$ws.Range("M122").Value = -23153816.5
$ws.Range("N122").Value = -8197755.25

# Row 131
$ws.Range("H131").Value = 2199.7778
$ws.Range("I131").Value = 2293.6
$ws.Range("J131").Value = 2178.4546
$ws.Range("K131").Value = 6880.799999999999
$ws.Range("L131").Value = 6535.3638
$ws.Range("M131").Value = -1840.799999999999
$ws.Range("N131").Value = -16615.3638

# Row 132
$ws.Range("H132").Value = 6773.173
$ws.Range("I132").Value = 3836.7307
$ws.Range("J132").Value = 9709.615
$ws.Range("K132").Value = 34530.5763
$ws.Range("L132").Value = 87386.535
$ws.Range("M132").Value = -32000.5763
$ws.Range("N132").Value = -92446.535

# Row 137
$ws.Range("H137").Value = 103907.1
$ws.Range("I137").Value = 3575.7144
$ws.Range("J137").Value = 338013.66
$ws.Range("K137").Value = 10727.1432
$ws.Range("L137").Value = 1014040.98
$ws.Range("M137").Value = -5627.143199999999
$ws.Range("N137").Value = -1024240.98

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = $null
$ws.Range("N44").Value = 0

# Row 106
$ws.Range("H106").Value = 26999
$ws.Range("J106").Value = 26999
$ws.Range("L106").Value = 26999
$ws.Range("N106").Value = -29523

# Row 132
$ws.Range("H132").Value = 5527.8335
$ws.Range("I132").Value = 3693.5715
$ws.Range("J132").Value = 11947.75
$ws.Range("K132").Value = 11080.7145
$ws.Range("L132").Value = 35843.25
$ws.Range("M132").Value = -8550.7145
$ws.Range("N132").Value = -40903.25

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 122
$ws.Range("H122").Value = 4359.3955
$ws.Range("I122").Value = 3428.724
$ws.Range("J122").Value = 6287.2144
$ws.Range("K122").Value = 10286.172
$ws.Range("L122").Value = 18861.6432
$ws.Range("M122").Value = -7836.172
$ws.Range("N122").Value = -23761.6432

# Row 132
$ws.Range("H132").Value = 7358797
$ws.Range("I132").Value = 11115135
$ws.Range("J132").Value = 9439.565000000001
$ws.Range("K132").Value = 33345405
$ws.Range("L132").Value = 28318.695
$ws.Range("M132").Value = -33342875
$ws.Range("N132").Value = -33378.695

# Row 136
$ws.Range("H136").Value = 14011.711
$ws.Range("I136").Value = 4031.5293
$ws.Range("K136").Value = 12094.5879
$ws.Range("M136").Value = -9544.5879

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 100
$ws.Range("H100").Value = 815.4167
$ws.Range("I100").Value = 989.6429000000001
$ws.Range("K100").Value = 1979.2858
$ws.Range("M100").Value = -1438.2858

# Row 107
$ws.Range("H107").Value = 11495330
$ws.Range("I107").Value = 590.94446
$ws.Range("J107").Value = 30304904
$ws.Range("K107").Value = 1772.83338
$ws.Range("L107").Value = 90914712
$ws.Range("M107").Value = 147.16662
$ws.Range("N107").Value = -90918552

# Row 113
$ws.Range("I113").Value = 1265
$ws.Range("J113").Value = 2823.2
$ws.Range("K113").Value = 3795
$ws.Range("L113").Value = 8469.599999999999
$ws.Range("M113").Value = -1625
$ws.Range("N113").Value = -12809.6

# Row 122
$ws.Range("H122").Value = 3709.7437
$ws.Range("I122").Value = 2128.7144
$ws.Range("K122").Value = 6386.1432
$ws.Range("M122").Value = -3936.1432

# Row 125
$ws.Range("H125").Value = 61715
$ws.Range("J125").Value = 61715
$ws.Range("L125").Value = 61715
$ws.Range("N125").Value = -71555

# Row 132
$ws.Range("H132").Value = 7901.4
$ws.Range("I132").Value = 9372.519
$ws.Range("K132").Value = 28117.557
$ws.Range("M132").Value = -25587.557

# Row 136
$ws.Range("H136").Value = 15389595
$ws.Range("I136").Value = 27030576
$ws.Range("K136").Value = 81091728
$ws.Range("M136").Value = -81089178
